$wb = $excel.ActiveWorkbook

# --- Sheet3 (OrganizationTestData): add new Organization test case block ---
$ws3 = $wb.Worksheets.Item("OrganizationTestData")

$ws3.Range("A5").Value = "TC_ID"
$ws3.Range("B5").Value = "TC_Name"
$ws3.Range("C5").Value = "Attribute"
$ws3.Range("D5").Value = "Data"

$ws3.Range("A6").Value = "TC_Contacts_02"
$ws3.Range("B6").Value = "Create Contact With Child Windows Details In Contact Information"
$ws3.Range("C6").Value = "OrganizationName"
$ws3.Range("D6").Value = "Amazon"

$ws3.Range("C7").Value = "Industry"
$ws3.Range("D7").Value = "Finance"

$ws3.Range("C8").Value = "Type"
$ws3.Range("D8").Value = "Investor"

$ws3.Range("C9").Value = "Rating"
$ws3.Range("D9").Value = "Acquired"

$ws3.Range("A6:A11").Merge()
$ws3.Range("B6:B8").Merge()

# --- Sheet2 (ContactsTestData): D6 "Shaktiman" -> "Potter" ---
$ws2 = $wb.Worksheets.Item("ContactsTestData")
$ws2.Range("D6").Value = "Potter"
